$d = $word.ActiveDocument
$ldq = [char]0x201C    # U+201C LEFT DOUBLE QUOTATION MARK “

# --- 1) Version history table: merge "Version 1" / "." / "6" -> "Version 1.6" ---
# (text itself does not change, Find/Replace with identical text normalizes the runs
#  the same way the real edit did; harmless no-op if already merged)
$d.Content.Find.Execute("Version 1.6", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version 1.6", 2) | Out-Null

# --- 2) Version history table: merge "Update after" / " an internal rewiew" ---
$d.Content.Find.Execute("Update after an internal rewiew", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Update after an internal rewiew", 2) | Out-Null

# --- 3) SRS_Home_1 paragraph: merge the split runs (no visible text change) ---
$d.Content.Find.Execute("The user can select a trip from a photo gallery of the top travel destinations. ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "The user can select a trip from a photo gallery of the top travel destinations. ", 2) | Out-Null

# --- 4) SRS_Home_2 paragraph: "The user can navigate..." -> "The registered user can navigate..." ---
$rngHome2 = $d.Content
$foundHome2 = $rngHome2.Find.Execute("The user can navigate to the booking", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if ($foundHome2) {
    # Insert " registered" right after the leading "The" (matches how the author typed the edit)
    $afterThe = $d.Range($rngHome2.Start + 3, $rngHome2.Start + 3)
    $afterThe.InsertAfter(" registered")
}

# --- 5) SRS_Admn_6 paragraph: prepend "When the admin deletes user, " ---
$findAdmn6 = "The program should display a message containing " + $ldq + "User added "
$replAdmn6 = "When the admin deletes user, the program should display a message containing " + $ldq + "User added "
$d.Content.Find.Execute($findAdmn6, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replAdmn6, 2) | Out-Null

# --- 6) Remove the whole "SRS_Admn_12:" requirement (two paragraphs) ---
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text
    if ($txt -match "^SRS_Admn_12:") {
        $next = $paras.Item($i + 1)
        $r = $d.Range($p.Range.Start, $next.Range.End)
        $r.Delete()
        break
    }
}

# --- 7) "SRS_BK_1:" label: merge the trailing "_1" / ":" runs (no visible text change) ---
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "SRS_BK_1:`r") {
        $p.Range.Find.Execute("_1:", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "_1:", 2) | Out-Null
        break
    }
}

# --- 8) SRS_BK_1 paragraph: "The user can book..." -> "Registered user can book..." ---
$d.Content.Find.Execute("The user can book any flight and a message appears", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Registered user can book any flight and a message appears", 2) | Out-Null
